# Applies cryptos price/volume/coin updates per commit "Updated symbol list on Mon Feb  6 19:31:49 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that should be stored as literal text (prices/percentages),
# avoiding Excel auto-converting numeric-looking strings into numbers/percentages.
function Set-TextValue($cell, [string]$val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "329.16"
Set-TextValue $ws.Range("E2") "0.29%"

# Row 3
Set-TextValue $ws.Range("D3") "44.36"
Set-TextValue $ws.Range("E3") "1.00%"

# Row 4
Set-TextValue $ws.Range("D4") "5.601"
Set-TextValue $ws.Range("E4") "3.56%"

# Row 5
Set-TextValue $ws.Range("D5") "0.08099"
Set-TextValue $ws.Range("E5") "0.11%"

# Row 6
Set-TextValue $ws.Range("D6") "2.014"
Set-TextValue $ws.Range("E6") "5.99%"

# Row 7
Set-TextValue $ws.Range("D7") "4.309"
Set-TextValue $ws.Range("E7") "0.08%"

# Row 8
Set-TextValue $ws.Range("D8") "0.9545"
Set-TextValue $ws.Range("E8") "1.40%"

# Row 9
Set-TextValue $ws.Range("E9") "-5.82%"

# Row 10
Set-TextValue $ws.Range("D10") "0.1191"
Set-TextValue $ws.Range("E10") "-1.41%"

# Row 11
Set-TextValue $ws.Range("D11") "0.1859"
Set-TextValue $ws.Range("E11") "-1.85%"

# Row 12
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D12") "0.09738"
Set-TextValue $ws.Range("E12") "2.70%"

# Row 13
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D13") "0.04596"
Set-TextValue $ws.Range("E13") "10.89%"

# Row 14
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D14") "0.1069"
Set-TextValue $ws.Range("E14") "-0.19%"

# Row 15
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D15") "0.001273"
Set-TextValue $ws.Range("E15") "-0.11%"

# Row 16
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws.Range("D16") "0.04218"
Set-TextValue $ws.Range("E16") "-3.49%"

# Row 17
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D17") "0.005894"
Set-TextValue $ws.Range("E17") "-3.09%"

# Row 18
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D18") "3.369"
Set-TextValue $ws.Range("E18") "-5.73%"

# Row 19
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue $ws.Range("D19") "0.3474"
Set-TextValue $ws.Range("E19") "-0.70%"

# Row 20
$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue $ws.Range("D20") "10.20"
Set-TextValue $ws.Range("E20") "19.81%"

# Row 21
Set-TextValue $ws.Range("E21") "3.67%"

# Row 22
Set-TextValue $ws.Range("D22") "0.2505"
Set-TextValue $ws.Range("E22") "-3.89%"

# Row 23
Set-TextValue $ws.Range("D23") "0.001248"
Set-TextValue $ws.Range("E23") "0.73%"

# Row 24
Set-TextValue $ws.Range("D24") "0.004324"
Set-TextValue $ws.Range("E24") "0.77%"

# Row 25
Set-TextValue $ws.Range("D25") "0.0001190"
Set-TextValue $ws.Range("E25") "-3.70%"

# Row 26
Set-TextValue $ws.Range("E26") "-0.97%"

# Row 38
Set-TextValue $ws.Range("E38") "0.36%"

# Row 39
Set-TextValue $ws.Range("D39") "0.05560"
Set-TextValue $ws.Range("E39") "2.57%"

# Row 40
Set-TextValue $ws.Range("D40") "0.007591"
Set-TextValue $ws.Range("E40") "-1.30%"

# Row 41
Set-TextValue $ws.Range("D41") "0.1409"
Set-TextValue $ws.Range("E41") "1.16%"

# Row 42
Set-TextValue $ws.Range("D42") "0.008080"
Set-TextValue $ws.Range("E42") "-17.25%"

# Row 43
Set-TextValue $ws.Range("D43") "0.002015"
Set-TextValue $ws.Range("E43") "-5.37%"

# Row 44
Set-TextValue $ws.Range("D44") "0.008402"
Set-TextValue $ws.Range("E44") "-15.04%"

# Row 45
Set-TextValue $ws.Range("D45") "0.00007186"
Set-TextValue $ws.Range("E45") "1.46%"

# Row 46
Set-TextValue $ws.Range("D46") "0.00000000750"
Set-TextValue $ws.Range("E46") "-0.72%"

# Row 47
Set-TextValue $ws.Range("D47") "0.004184"
Set-TextValue $ws.Range("E47") "17.89%"

# Row 48
Set-TextValue $ws.Range("D48") "0.002270"
Set-TextValue $ws.Range("E48") "-0.71%"

# Row 49
Set-TextValue $ws.Range("D49") "0.00002100"
Set-TextValue $ws.Range("E49") "-0.72%"

# Row 50
Set-TextValue $ws.Range("D50") "0.0002000"
Set-TextValue $ws.Range("E50") "-0.72%"
